$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear PAN for Investor 1 (row 2) - PAN is nil
$ws.Range("D2").Clear()

# Clear Investing Entity (name) and PAN for Investor 3 (row 4) - name/PAN are nil
$ws.Range("B4").Clear()
$ws.Range("D4").Clear()

# Update selected cell to D4
$ws.Range("D4").Select()
